# Auto-generated Excel COM-interop script to update cryptos.xlsx price data
# Sets B/C/D/E cell text values for rows 2-51 as plain text, preserving
# number-formatted strings (e.g. trailing zeros, European thousand separators)
# exactly as authored, and avoiding Excel auto-converting numeric-looking
# strings into floating point numbers or changing cell styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($address, $value) {
    $cell = $ws.Range($address)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell "D2" "64.074.31"
Set-TextCell "E2" "  -1.94%  "
Set-TextCell "D3" "3.350.01"
Set-TextCell "E3" "  -2.66%  "
Set-TextCell "D4" "0.999"
Set-TextCell "E4" "  -0.15%  "
Set-TextCell "D5" "551.07"
Set-TextCell "E5" "  -0.22%  "
Set-TextCell "D6" "173.37"
Set-TextCell "E6" "  -2.60%  "
Set-TextCell "D7" "0.614"
Set-TextCell "E7" "  -3.42%  "
Set-TextCell "D8" "3.341.57"
Set-TextCell "E8" "  -2.78%  "
Set-TextCell "E9" "  -0.01%  "
Set-TextCell "D10" "0.614"
Set-TextCell "E10" "  -3.10%  "
Set-TextCell "D11" "0.153"
Set-TextCell "E11" "  +0.53%  "
Set-TextCell "D12" "53.98"
Set-TextCell "E12" "  +0.25%  "
Set-TextCell "D13" "0.0000266"
Set-TextCell "E13" "  -1.63%  "
Set-TextCell "D14" "8.91"
Set-TextCell "E14" "  -3.27%  "
Set-TextCell "D15" "3.729.30"
Set-TextCell "E15" "  -6.74%  "
Set-TextCell "D16" "18.22"
Set-TextCell "E16" "  -1.43%  "
Set-TextCell "D17" "3.353.32"
Set-TextCell "E17" "  -2.70%  "
Set-TextCell "E18" "  -2.91%  "
Set-TextCell "D19" "11.75"
Set-TextCell "E19" "  -1.09%  "
Set-TextCell "D20" "63.919.86"
Set-TextCell "E20" "  -2.25%  "
Set-TextCell "D21" "0.980"
Set-TextCell "E21" "  -0.80%  "
Set-TextCell "D22" "411.00"
Set-TextCell "E22" "  -1.02%  "
Set-TextCell "D23" "4.05"
Set-TextCell "E23" "  +0.30%  "
Set-TextCell "D24" "4.40"
Set-TextCell "E24" "  +6.52%  "
Set-TextCell "D25" "13.88"
Set-TextCell "E25" "  +10.87%  "
Set-TextCell "D26" "83.09"
Set-TextCell "E26" "  -3.64%  "
Set-TextCell "D27" "10.58"
Set-TextCell "E27" "  -1.90%  "
Set-TextCell "D28" "2.74"
Set-TextCell "E28" "  -4.29%  "
Set-TextCell "D29" "8.65"
Set-TextCell "E29" "  -4.55%  "
Set-TextCell "D30" "29.18"
Set-TextCell "E30" "  -2.82%  "
Set-TextCell "D31" "6.41"
Set-TextCell "E31" "  -2.18%  "
Set-TextCell "D32" "581.46"
Set-TextCell "E32" "  -4.26%  "
Set-TextCell "D33" "11.38"
Set-TextCell "E33" "  -2.89%  "
Set-TextCell "D34" "0.107"
Set-TextCell "E34" "  -2.40%  "
Set-TextCell "D35" "58.07"
Set-TextCell "E35" "  -1.59%  "
Set-TextCell "D36" "0.148"
Set-TextCell "E36" "  +2.58%  "
Set-TextCell "E37" "  -0.02%  "
Set-TextCell "D38" "35.28"
Set-TextCell "E38" "  -5.59%  "
Set-TextCell "B39" "PEPE"
Set-TextCell "C39" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell "D39" "0.0₃0741"
Set-TextCell "E39" "  -5.72%  "
Set-TextCell "B40" "Stacks"
Set-TextCell "C40" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell "D40" "3.40"
Set-TextCell "E40" "  +2.68%  "
Set-TextCell "D41" "0.368"
Set-TextCell "E41" "  -2.92%  "
Set-TextCell "D42" "3.145.10"
Set-TextCell "E42" "  -2.42%  "
Set-TextCell "D43" "0.997"
Set-TextCell "E43" "  -0.38%  "
Set-TextCell "D44" "2.81"
Set-TextCell "E44" "  +0.27%  "
Set-TextCell "D45" "3.28"
Set-TextCell "E45" "  +1.06%  "
Set-TextCell "D46" "0.0401"
Set-TextCell "E46" "  -2.81%  "
Set-TextCell "D47" "2.43"
Set-TextCell "E47" "  -4.66%  "
Set-TextCell "D48" "2.61"
Set-TextCell "E48" "  -3.68%  "
Set-TextCell "E49" "  -3.26%  "
Set-TextCell "D50" "132.74"
Set-TextCell "E50" "  -3.61%  "
Set-TextCell "D51" "8.09"
Set-TextCell "E51" "  -3.82%  "
